$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 507, shifting existing rows 507:608 down to 508:609
$ws.Rows.Item(507).Insert()

# Populate the newly inserted row 507 with the new weekly record
$ws.Cells.Item(507, 1).Value = 1
$ws.Cells.Item(507, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(507, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(507, 4).Value = 44644
$ws.Cells.Item(507, 5).Value = 15
$ws.Cells.Item(507, 6).Value = 100112024
$ws.Cells.Item(507, 7).Value = "Choclo"
$ws.Cells.Item(507, 8).Value = "Sin especificar"
$ws.Cells.Item(507, 9).Value = "Primera"
$ws.Cells.Item(507, 10).Value = 70
$ws.Cells.Item(507, 11).Value = 42000
$ws.Cells.Item(507, 12).Value = 45000
$ws.Cells.Item(507, 13).Value = 43500
$ws.Cells.Item(507, 14).Value = "$/saco 100 unidades"
$ws.Cells.Item(507, 15).Value = "Valle de Camiña"
$ws.Cells.Item(507, 16).Value = 435
$ws.Cells.Item(507, 17).Value = 100
$ws.Cells.Item(507, 18).Value = "Hortaliza"
